$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '22.021.41'
$ws.Cells.Item(2, 5).Value = '  -1.92%  '

$ws.Cells.Item(3, 4).Value = '1.554.50'

$ws.Cells.Item(4, 5).Value = '  +0.03%  '

$ws.Cells.Item(5, 5).Value = '  +0.02%  '

$cell = $ws.Cells.Item(6, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '286.66'
$cell.Style = $origStyle
$ws.Cells.Item(6, 5).Value = '  -0.47%  '

$cell = $ws.Cells.Item(7, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.3756'
$cell.Style = $origStyle
$ws.Cells.Item(7, 5).Value = '  +0.92%  '

$cell = $ws.Cells.Item(8, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.3238'
$cell.Style = $origStyle
$ws.Cells.Item(8, 5).Value = '  -2.46%  '

$ws.Cells.Item(9, 2).Value = 'OKB'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$cell = $ws.Cells.Item(9, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '41.37'
$cell.Style = $origStyle
$ws.Cells.Item(9, 5).Value = '  -13.01%  '

$ws.Cells.Item(10, 2).Value = 'Polygon'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell = $ws.Cells.Item(10, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.126'
$cell.Style = $origStyle
$ws.Cells.Item(10, 5).Value = '  -2.55%  '

$cell = $ws.Cells.Item(11, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.07297'
$cell.Style = $origStyle
$ws.Cells.Item(11, 5).Value = '  -3.00%  '

$cell = $ws.Cells.Item(12, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = $origStyle
$ws.Cells.Item(12, 5).Value = '  +0.02%  '

$cell = $ws.Cells.Item(13, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '19.54'
$cell.Style = $origStyle
$ws.Cells.Item(13, 5).Value = '  -6.02%  '

$cell = $ws.Cells.Item(14, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.705'
$cell.Style = $origStyle
$ws.Cells.Item(14, 5).Value = '  -3.87%  '

$cell = $ws.Cells.Item(15, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.854'
$cell.Style = $origStyle
$ws.Cells.Item(15, 5).Value = '  -1.07%  '

$ws.Cells.Item(16, 4).Value = '1.552.06'
$ws.Cells.Item(16, 5).Value = '  -0.91%  '

$cell = $ws.Cells.Item(17, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.00001080'
$cell.Style = $origStyle
$ws.Cells.Item(17, 5).Value = '  -3.38%  '

$cell = $ws.Cells.Item(18, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.06640'
$cell.Style = $origStyle
$ws.Cells.Item(18, 5).Value = '  -1.29%  '

$cell = $ws.Cells.Item(19, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '85.05'
$cell.Style = $origStyle
$ws.Cells.Item(19, 5).Value = '  -3.72%  '

$cell = $ws.Cells.Item(20, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.446'
$cell.Style = $origStyle
$ws.Cells.Item(20, 5).Value = '  +0.83%  '

$cell = $ws.Cells.Item(21, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = $origStyle
$ws.Cells.Item(21, 5).Value = '  +0.04%  '

$cell = $ws.Cells.Item(22, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '15.97'
$cell.Style = $origStyle
$ws.Cells.Item(22, 5).Value = '  -3.27%  '

$cell = $ws.Cells.Item(23, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.57'
$cell.Style = $origStyle
$ws.Cells.Item(23, 5).Value = '  -3.71%  '

$ws.Cells.Item(24, 4).Value = '22.051.88'
$ws.Cells.Item(24, 5).Value = '  -1.77%  '

$cell = $ws.Cells.Item(25, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.244'
$cell.Style = $origStyle
$ws.Cells.Item(25, 5).Value = '  -6.35%  '

$cell = $ws.Cells.Item(26, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.523'
$cell.Style = $origStyle
$ws.Cells.Item(26, 5).Value = '  -3.89%  '

$cell = $ws.Cells.Item(27, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '149.51'
$cell.Style = $origStyle
$ws.Cells.Item(27, 5).Value = '  -0.73%  '

$cell = $ws.Cells.Item(28, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '18.84'
$cell.Style = $origStyle
$ws.Cells.Item(28, 5).Value = '  -4.09%  '

$cell = $ws.Cells.Item(29, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.823'
$cell.Style = $origStyle
$ws.Cells.Item(29, 5).Value = '  -2.71%  '

$ws.Cells.Item(30, 4).Value = '1.731.94'
$ws.Cells.Item(30, 5).Value = '  -0.82%  '

$cell = $ws.Cells.Item(31, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '119.89'
$cell.Style = $origStyle
$ws.Cells.Item(31, 5).Value = '  -4.26%  '

$cell = $ws.Cells.Item(32, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.117'
$cell.Style = $origStyle
$ws.Cells.Item(32, 5).Value = '  +1.85%  '

$cell = $ws.Cells.Item(33, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.923'
$cell.Style = $origStyle
$ws.Cells.Item(33, 5).Value = '  -2.75%  '

$cell = $ws.Cells.Item(34, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.249'
$cell.Style = $origStyle
$ws.Cells.Item(34, 5).Value = '  -5.95%  '

$cell = $ws.Cells.Item(35, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.08108'
$cell.Style = $origStyle
$ws.Cells.Item(35, 5).Value = '  -2.74%  '

$cell = $ws.Cells.Item(36, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.633'
$cell.Style = $origStyle
$ws.Cells.Item(36, 5).Value = '  -17.94%  '

$ws.Cells.Item(37, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Cells.Item(37, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.220'
$cell.Style = $origStyle
$ws.Cells.Item(37, 5).Value = '  -2.40%  '

$ws.Cells.Item(38, 2).Value = 'VeChain'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Cells.Item(38, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.02283'
$cell.Style = $origStyle
$ws.Cells.Item(38, 5).Value = '  -7.03%  '

$cell = $ws.Cells.Item(39, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.06118'
$cell.Style = $origStyle
$ws.Cells.Item(39, 5).Value = '  -4.36%  '

$cell = $ws.Cells.Item(40, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.2108'
$cell.Style = $origStyle
$ws.Cells.Item(40, 5).Value = '  -5.60%  '

$ws.Cells.Item(41, 5).Value = '  -7.34%  '

$cell = $ws.Cells.Item(42, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.90'
$cell.Style = $origStyle
$ws.Cells.Item(42, 5).Value = '  -4.46%  '

$ws.Cells.Item(43, 5).Value = '  +0.05%  '

$cell = $ws.Cells.Item(44, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5936'
$cell.Style = $origStyle
$ws.Cells.Item(44, 5).Value = '  -5.32%  '

$cell = $ws.Cells.Item(45, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '13.58'
$cell.Style = $origStyle
$ws.Cells.Item(45, 5).Value = '  -2.71%  '

$ws.Cells.Item(46, 5).Value = '  -1.29%  '

$cell = $ws.Cells.Item(47, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5744'
$cell.Style = $origStyle
$ws.Cells.Item(47, 5).Value = '  -5.71%  '

$ws.Cells.Item(48, 5).Value = '  -5.22%  '

$cell = $ws.Cells.Item(49, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '119.91'
$cell.Style = $origStyle
$ws.Cells.Item(49, 5).Value = '  -4.09%  '

$cell = $ws.Cells.Item(50, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.154'
$cell.Style = $origStyle
$ws.Cells.Item(50, 5).Value = '  -4.59%  '

$cell = $ws.Cells.Item(51, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.06929'
$cell.Style = $origStyle
$ws.Cells.Item(51, 5).Value = '  -3.80%  '
